$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.245.32"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.593.71"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.90"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.90"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0853"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.817.10"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.597.01"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.504"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.58"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.208.76"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.23"
$ws.Range("E18").Value = "  +7.07%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.58"
$ws.Range("E20").Value = "  +4.28%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.24"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.14"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.91"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.57"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.97"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.32"
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.462.29"
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.95"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.565"
$ws.Range("E37").Value = "  -4.05%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.818"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.76"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  +1.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.926"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.729.01"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.754"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.46"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.55"
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0500"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0945"
$ws.Range("E51").Value = "  -2.67%  "
